$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "'310307356"
$ws.Range("G2").Value = "'320069738"

$ws.Range("G7").Select()
